# Refresh the crypto price/volume table (coinranking.com snapshot).
# Source cells in columns B:E are stored as plain text (never numeric), so
# every new value below is written with a leading apostrophe - Excel's
# "force text" entry convention - so values such as "27.621.54", "1.0000"
# or "0.9991" are kept as literal text instead of being auto-coerced into
# numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
}

Set-TextValue "D2" "27.621.54"
Set-TextValue "E2" "  -1.60%  "

Set-TextValue "D3" "1.878.22"
Set-TextValue "E3" "  -1.42%  "

Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  -0.18%  "

Set-TextValue "D5" "329.97"
Set-TextValue "E5" "  +0.71%  "

Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  -0.20%  "

Set-TextValue "D7" "0.4713"
Set-TextValue "E7" "  +1.99%  "

Set-TextValue "D8" "0.3988"
Set-TextValue "E8" "  +0.29%  "

Set-TextValue "D9" "47.34"
Set-TextValue "E9" "  -8.91%  "

Set-TextValue "D10" "0.08069"
Set-TextValue "E10" "  -3.34%  "

Set-TextValue "D11" "1.027"
Set-TextValue "E11" "  -1.29%  "

Set-TextValue "D12" "21.88"
Set-TextValue "E12" "  +0.05%  "

Set-TextValue "D13" "1.862.82"
Set-TextValue "E13" "  -2.35%  "

Set-TextValue "D14" "5.964"
Set-TextValue "E14" "  -0.99%  "

Set-TextValue "D15" "7.215"
Set-TextValue "E15" "  -2.19%  "

Set-TextValue "D16" "1.003"
Set-TextValue "E16" "  +0.05%  "

Set-TextValue "D17" "87.02"
Set-TextValue "E17" "  -2.39%  "

Set-TextValue "D18" "0.00001042"
Set-TextValue "E18" "  -2.16%  "

Set-TextValue "D19" "0.06575"
Set-TextValue "E19" "  -0.24%  "

Set-TextValue "D20" "17.35"
Set-TextValue "E20" "  -2.60%  "

Set-TextValue "D21" "0.9991"
Set-TextValue "E21" "  -0.14%  "

Set-TextValue "D22" "5.524"
Set-TextValue "E22" "  -3.28%  "

Set-TextValue "D23" "27.638.28"
Set-TextValue "E23" "  -1.57%  "

Set-TextValue "D24" "11.02"
Set-TextValue "E24" "  -1.19%  "

Set-TextValue "D25" "2.301"
Set-TextValue "E25" "  -0.51%  "

Set-TextValue "D26" "2.079.89"
Set-TextValue "E26" "  -2.41%  "

Set-TextValue "B27" "EthereumClassic"
Set-TextValue "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D27" "20.36"
Set-TextValue "E27" "  +1.96%  "

Set-TextValue "B28" "Monero"
Set-TextValue "C28" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D28" "154.21"
Set-TextValue "E28" "  +0.47%  "

Set-TextValue "D29" "2.102"
Set-TextValue "E29" "  -0.96%  "

Set-TextValue "D30" "5.558"
Set-TextValue "E30" "  -2.83%  "

Set-TextValue "D31" "122.64"
Set-TextValue "E31" "  -0.59%  "

Set-TextValue "D32" "0.09502"
Set-TextValue "E32" "  -0.92%  "

Set-TextValue "D33" "0.9582"
Set-TextValue "E33" "  -1.12%  "

Set-TextValue "D34" "1.478"
Set-TextValue "E34" "  +0.41%  "

Set-TextValue "D35" "3.602"
Set-TextValue "E35" "  -0.52%  "

Set-TextValue "D36" "5.323"
Set-TextValue "E36" "  -3.36%  "

Set-TextValue "D37" "0.06118"
Set-TextValue "E37" "  -0.18%  "

Set-TextValue "D38" "0.02259"
Set-TextValue "E38" "  -1.06%  "

Set-TextValue "D39" "1.221"
Set-TextValue "E39" "  -3.49%  "

Set-TextValue "D40" "8.280"
Set-TextValue "E40" "  -4.76%  "

Set-TextValue "D41" "0.6012"
Set-TextValue "E41" "  -1.88%  "

Set-TextValue "D42" "1.0000"
Set-TextValue "E42" "  -0.15%  "

Set-TextValue "E43" "  +0.11%  "

Set-TextValue "D44" "10.38"
Set-TextValue "E44" "  -4.19%  "

Set-TextValue "B45" "Decentraland"
Set-TextValue "C45" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D45" "0.5708"
Set-TextValue "E45" "  -2.28%  "

Set-TextValue "B46" "WEMIXTOKEN"
Set-TextValue "C46" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D46" "1.247"
Set-TextValue "E46" "  -4.46%  "

Set-TextValue "D47" "12.20"
Set-TextValue "E47" "  -4.71%  "

Set-TextValue "D48" "3.411"
Set-TextValue "E48" "  -0.67%  "

Set-TextValue "D49" "1.943"
Set-TextValue "E49" "  -2.98%  "

Set-TextValue "D50" "0.06818"
Set-TextValue "E50" "  -1.26%  "

Set-TextValue "D51" "110.10"
Set-TextValue "E51" "  -0.33%  "
